# Update countries & provincias Spain
# This script refreshes the COVID-19 "Pais" dashboard data:
#  - updates the "last updated" timestamp
#  - updates numeric statistics (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes hoy, Muertes) for the rows whose
#    source data changed
#  - because the underlying ranking shifted slightly, a handful of rows now
#    show a different country name (the row "slot" stayed the same, but the
#    country associated with that slot changed)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Header / timestamp -------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 17 de Octubre de 2020 a las 18:46"

# Row data: Row, Country, CasosTotales, NuevosCasos, CasosActivos, Recuperados, CasosCriticos, MuertesHoy, Muertes
$rows = @(
    @(4,   "Estados Unidos",  8313116, 24838, 5405763, 2683407, 0, 302, 223946),
    @(5,   "India",           7478924, 48289, 6577729,  787582, 0, 581, 113613),
    @(6,   "Brasil",          5205686,  4116, 4619560,  432768, 0, 129, 153358),
    @(21,  "Alemania",         359558,  2766,  290000,   59713, 0,   9,   9845),
    @(24,  "Turquia",          345678,  1723,  302499,   33955, 0,  71,   9224),
    @(27,  "Israel",           302730,   834,  265348,   35215, 0,  26,   2167),
    @(31,  "Canada",           196192,  2086,  165417,   21029, 0,  24,   9746),
    @(35,  "Chequia",          164422,  4310,   68896,   94188, 0,  55,   1338),
    @(77,  "Jordania",          36053,  1505,    6773,   28950, 0,  20,    330),
    @(78,  "Serbia",            35946,   227,   31536,    3636, 0,   2,    774),
    @(79,  "Dinamarca",         34941,   500,   28917,    5345, 0,   2,    679),
    @(80,  "Birmania",          34875,  1387,   16370,   17667, 0,  39,    838),
    @(81,  "Tunez",             34790,     0,    5032,   29246, 0,   0,    512),
    @(88,  "Grecia",            24932,   482,    9989,   14443, 0,  10,    500),
    @(89,  "Croacia",           24761,  1096,   19562,    4844, 0,  10,    355),
    @(99,  "Montenegro",        15427,    74,   10768,    4428, 0,   3,    231),
    @(100, "Senegal",           15392,    24,   13756,    1319, 0,   0,    317),
    @(108, "Mozambique",        10707,    95,    8272,    2361, 0,   1,     74),
    @(109, "Luxemburgo",        10646,   175,    8468,    2045, 0,   0,    133),
    @(159, "Sierra Leona",       2327,     2,    1753,     501, 0,   0,     73),
    @(165, "Liberia",            1377,     0,    1268,      27, 0,   0,     82),
    @(209, "Granada",              27,     2,      24,       3, 0,   0,      0),
    @(210, "Nueva Caledonia",      27,     0,      27,       0, 0,   0,      0),
    @(216, "Montserrat",           13,     0,      12,       0, 0,   0,      1),
    @(217, "Islas Malvinas",       13,     0,      13,       0, 0,   0,      0)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
